$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update method_type labels (column C) and colours (column E) for rows 12-14
$ws.Range("C12").Value = "Empirical - other"
$ws.Range("E12").Value = "#33a02c"

$ws.Range("C13").Value = "Empirical - social primary"
$ws.Range("E13").Value = "#b2df8a"

$ws.Range("E14").Value = "#1f78b4"

# Resize column A to fit new content (matches explicit width seen in target)
$ws.Columns.Item(1).ColumnWidth = 16.3

# Update view: scroll back to top-left and change selection
$ws.Range("C14").Select()

$wb.Save()
